# This script applies a row-rotation update to rows 5-13 of the active sheet.
# Each of these rows describes one "artfynd" (species observation) record.
# The records' identifying/describing columns (A, B, E, F, G, H, Q, R) get
# redistributed among the rows (a permutation), with column B (taxon sort
# order) additionally shifted by +14 except for the "bird" record, and the
# two special comment cells (M, AC) belonging to the bird record move along
# with the rest of that record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that carry the per-record data that gets redistributed.
$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")

# Snapshot the current ("before") values of every relevant cell, for rows 5-13,
# before any writes happen (writes happen later, so this is safe either way,
# but we snapshot explicitly to make the row permutation correct/explicit).
$snapshot = @{}
foreach ($r in 5..13) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}
$snapshotM10 = $ws.Range("M10").Value2
$snapshotAC10 = $ws.Range("AC10").Value2

# Target row <- source row mapping: the data that used to live in row
# $mapping[$target] now belongs in row $target.
$mapping = @{
    5  = 6
    6  = 13
    7  = 5
    8  = 7
    9  = 12
    10 = 11
    11 = 9
    12 = 10
    13 = 8
}

foreach ($target in $mapping.Keys) {
    $source = $mapping[$target]
    $src = $snapshot[$source]

    $ws.Range("A$target").Value = $src["A"]
    if ($source -eq 10) {
        # The bird record keeps its original taxon sort order value.
        $ws.Range("B$target").Value = $src["B"]
    } else {
        $ws.Range("B$target").Value = $src["B"] + 14
    }
    $ws.Range("E$target").Value = $src["E"]
    $ws.Range("F$target").Value = $src["F"]
    $ws.Range("G$target").Value = $src["G"]
    $ws.Range("H$target").Value = $src["H"]
    $ws.Range("Q$target").Value = $src["Q"]
    $ws.Range("R$target").Value = $src["R"]
}

# Move the "activity"/"public comment" notes tied to the bird record
# (previously in row 10) along with the rest of its data, now in row 12.
$ws.Range("M10").Value = $null
$ws.Range("AC10").Value = $null
$ws.Range("M12").Value = $snapshotM10
$ws.Range("AC12").Value = $snapshotAC10
